$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 is no longer the last row, so its date cell should use the same
# number format as the rest of the date column ("YYYY-MM-DD HH:MM:SS",
# style index 2) instead of the "last row" format ("YYYY-MM-DD", style
# index 3).
$ws.Range("A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new last row (15) with the values from the diff, using the
# "last row" date format that row 14 previously had.
$ws.Range("A15").Value = 44526
$ws.Range("A15").NumberFormat = "YYYY-MM-DD"
$ws.Range("B15").Value = 551.5999999999985
